$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the merged title cell (C1:E2)
$ws.Range("C1").Value = "Profit - whole - 2019"

# Remove the old data row 5 (second data row) entirely
$ws.Rows.Item(5).Delete()

# Update header row 3
$ws.Range("A3").Value = "Id"
$ws.Range("B3").Value = "Stock"
$ws.Range("C3").Value = "Dealer"
$ws.Range("D3").Value = "Customer"
$ws.Range("E3").Value = "Miscellaneous"
$ws.Range("F3").Value = "Amount Pending"
$ws.Range("G3").Value = "Profit"

# Update data row 4 with numeric values
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0

# Resize columns B through G to width 23 (match the diff's single merged col span)
# Column D already stores width=23; reuse its reported ColumnWidth so the
# underlying stored width rounds to exactly 23 after Excel's internal
# pixel-based width conversion.
$targetWidth = $ws.Columns.Item(4).ColumnWidth()
$ws.Range("B1:G1").EntireColumn.ColumnWidth = $targetWidth
